$wb = $excel.ActiveWorkbook
$wsP = $wb.Worksheets.Item("parallelRuns")

# Rename header labels: "nThreads" -> "nProcesses" (B1 merged header, D11 table header)
$wsP.Cells.Item(1,2).Value = "nProcesses"
$wsP.Cells.Item(11,4).Value = "nProcesses"

# Add new "command" header in G11
$wsP.Cells.Item(11,7).Value = "command"
$wsP.Cells.Item(11,7).Font.Bold = $true

# Update run_name formula (col A) and add command formula (col G) for every data row
for ($r = 12; $r -le 71; $r++) {
    $wsP.Cells.Item($r,1).Formula = '="parallel_particles"&B' + $r + '&"_processes"&D' + $r + '&"_run"&E' + $r
    $wsP.Cells.Item($r,7).Formula = '="srun python parallel_mpi_pso_algorithm.py "&CHAR(34)&CHAR(34) &" 50 "&B' + $r + '&" 0.25 0.2 0.7 "&A' + $r + '&" "&D' + $r
}

# Re-apply the autofilter over the full data range with two column criteria
$wsP.AutoFilterMode = $false
$rng = $wsP.Range("A11:E71")
$rng.AutoFilter(2, @("140"))
$rng.AutoFilter(5, @("1"))

# Fix the hidden FilterDatabase defined name so it also spans the full range
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name() -eq "parallelRuns!_FilterDatabase") {
        $nm.RefersTo = "=parallelRuns!`$A`$11:`$E`$71"
    }
}

# Hide all rows that don't match the filter (nParticles=140 and runCount=1), i.e. everything
# except rows 27-31
for ($r = 12; $r -le 71; $r++) {
    if ($r -lt 27 -or $r -gt 31) {
        $wsP.Rows.Item($r).Hidden = $true
    }
}

# Autofit the new command column and select the new active cell
$wsP.Columns.Item(7).AutoFit()
$wsP.Activate()
$wsP.Range("C75").Select()

Write-Host "done"
